# Applies the cryptos list update described in the commit:
# "Updated cryptos list on Thu Oct 17 20:19:39 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.015.78'
$ws.Range("E2").Value = '  -1.26%  '

$ws.Range("D3").Value = '2.600.73'
$ws.Range("E3").Value = '  -0.80%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = "'589.78"
$ws.Range("E5").Value = '  -2.50%  '

$ws.Range("D6").Value = "'149.89"
$ws.Range("E6").Value = '  -3.00%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").Value = "'0.546"
$ws.Range("E8").Value = '  -1.20%  '

$ws.Range("D9").Value = '2.598.12'
$ws.Range("E9").Value = '  -0.89%  '

$ws.Range("D10").Value = "'0.125"
$ws.Range("E10").Value = '  -2.35%  '

$ws.Range("E11").Value = '  -0.06%  '

$ws.Range("D12").Value = "'5.16"
$ws.Range("E12").Value = '  -1.71%  '

$ws.Range("E13").Value = '  -3.18%  '

$ws.Range("D14").Value = "'27.26"
$ws.Range("E14").Value = '  -1.88%  '

$ws.Range("D15").Value = '3.070.59'
$ws.Range("E15").Value = '  -1.15%  '

$ws.Range("D16").Value = "'0.0000179"
$ws.Range("E16").Value = '  -4.75%  '

$ws.Range("D17").Value = '67.006.41'
$ws.Range("E17").Value = '  -1.26%  '

$ws.Range("D18").Value = '2.600.18'
$ws.Range("E18").Value = '  -0.62%  '

$ws.Range("D19").Value = "'363.65"
$ws.Range("E19").Value = '  -0.82%  '

$ws.Range("D20").Value = "'11.03"
$ws.Range("E20").Value = '  -1.62%  '

$ws.Range("D21").Value = "'7.32"
$ws.Range("E21").Value = '  -4.65%  '

$ws.Range("D22").Value = "'4.29"
$ws.Range("E22").Value = '  -0.46%  '

$ws.Range("D23").Value = "'4.83"
$ws.Range("E23").Value = '  -2.49%  '

$ws.Range("D24").Value = "'2.05"
$ws.Range("E24").Value = '  -0.92%  '

$ws.Range("D25").Value = "'72.89"
$ws.Range("E25").Value = '  +3.57%  '

$ws.Range("E26").Value = '  +0.03%  '

$ws.Range("D27").Value = "'9.92"
$ws.Range("E27").Value = '  +0.09%  '

$ws.Range("D28").Value = '2.728.86'
$ws.Range("E28").Value = '  -0.64%  '

$ws.Range("D29").Value = "'581.49"
$ws.Range("E29").Value = '  -0.07%  '

$ws.Range("E30").Value = '  +0.08%  '

$ws.Range("D31").Value = '0.0₃0981'
$ws.Range("E31").Value = '  -6.41%  '

$ws.Range("E32").Value = '  -5.24%  '

$ws.Range("D33").Value = "'7.63"
$ws.Range("E33").Value = '  -3.52%  '

$ws.Range("E34").Value = '  -3.77%  '

$ws.Range("E35").Value = '  -0.05%  '

$ws.Range("D36").Value = "'0.125"
$ws.Range("E36").Value = '  -4.71%  '

$ws.Range("E37").Value = '  -3.02%  '

$ws.Range("D38").Value = "'156.04"
$ws.Range("E38").Value = '  -1.02%  '

$ws.Range("D39").Value = "'18.96"
$ws.Range("E39").Value = '  -2.47%  '

$ws.Range("E40").Value = '  -1.50%  '

$ws.Range("D41").Value = "'1.85"
$ws.Range("E41").Value = '  -0.42%  '

$ws.Range("D42").Value = "'5.19"
$ws.Range("E42").Value = '  -3.45%  '

$ws.Range("D43").Value = "'17.09"
$ws.Range("E43").Value = '  +3.93%  '

$ws.Range("D44").Value = "'2.51"
$ws.Range("E44").Value = '  -4.23%  '

$ws.Range("D46").Value = "'152.55"
$ws.Range("E46").Value = '  -2.94%  '

$ws.Range("E47").Value = '  -1.89%  '

$ws.Range("E48").Value = '  -1.55%  '

$ws.Range("D49").Value = "'1.68"
$ws.Range("E49").Value = '  -2.84%  '

$ws.Range("D50").Value = "'0.0777"
$ws.Range("E50").Value = '  -1.59%  '

$ws.Range("D51").Value = "'21.44"
$ws.Range("E51").Value = '  +2.04%  '
